$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price values in column D are stored as plain text in this workbook
# (e.g. "42.607.44", "0.999"). When a replacement value parses as a normal
# floating point number (single decimal point) Excel would silently convert
# the cell to a numeric value (dropping trailing zeros, changing the cell
# type). Force those specific cells to the Text number format first so the
# literal text is preserved exactly, matching the source data.

$ws.Range("D2").Value = '42.607.44'
$ws.Range("E2").Value = '  -1.43%  '

$ws.Range("D3").Value = '2.532.53'
$ws.Range("E3").Value = '  -1.77%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.90'
$ws.Range("E5").Value = '  -2.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.64'
$ws.Range("E6").Value = '  +3.83%  '

$ws.Range("E7").Value = '  -1.14%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -2.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.89'
$ws.Range("E10").Value = '  +0.91%  '

$ws.Range("E11").Value = '  -1.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.33'
$ws.Range("E12").Value = '  -1.70%  '

$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").Value = '2.932.56'
$ws.Range("E14").Value = '  -1.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.95'
$ws.Range("E15").Value = '  +5.28%  '

$ws.Range("D16").Value = '2.552.97'
$ws.Range("E16").Value = '  -2.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.810'
$ws.Range("E17").Value = '  -4.04%  '

$ws.Range("D18").Value = '42.610.35'
$ws.Range("E18").Value = '  -1.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.75'
$ws.Range("E19").Value = '  -1.36%  '

$ws.Range("E20").Value = '  -1.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.26'
$ws.Range("E21").Value = '  -2.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.44'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.71'
$ws.Range("E23").Value = '  -4.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.89'
$ws.Range("E24").Value = '  -2.51%  '

$ws.Range("E25").Value = '  -1.26%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.02'
$ws.Range("E27").Value = '  -4.00%  '

$ws.Range("E28").Value = '  -3.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.15'
$ws.Range("E29").Value = '  -2.73%  '

$ws.Range("E30").Value = '  -1.44%  '

$ws.Range("E31").Value = '  -0.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.33'
$ws.Range("E32").Value = '  +0.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.73'
$ws.Range("E33").Value = '  +11.71%  '

$ws.Range("E35").Value = '  -2.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.45'
$ws.Range("E36").Value = '  -2.64%  '

$ws.Range("E37").Value = '  -4.93%  '

$ws.Range("E38").Value = '  -6.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.112'
$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("E40").Value = '  +0.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.29'
$ws.Range("E41").Value = '  +8.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.15'
$ws.Range("E42").Value = '  -1.81%  '

$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.28'
$ws.Range("E44").Value = '  +1.08%  '

$ws.Range("E45").Value = '  -1.81%  '

$ws.Range("D46").Value = '1.967.22'
$ws.Range("E46").Value = '  -1.63%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.88'
$ws.Range("E47").Value = '  -0.54%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '80.82'
$ws.Range("E48").Value = '  -2.80%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.729.14'
$ws.Range("E49").Value = '  -3.30%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.192'
$ws.Range("E50").Value = '  -1.07%  '

$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.853'
$ws.Range("E51").Value = '  +9.97%  '
